$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20 of data, following the pattern of the existing rows.
# Copy the formatting (style) of the row above for column A (date cell)
# then set the new value.
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").Value = 45986

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.46481303148316
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 2.509429409292352
